$wb = $excel.ActiveWorkbook

# "Generate Report for Handback" — the handback for f1f0eb99-5294-4e84-87c1-4da7274f6c19
# has been processed (it is stale vs. the latest commit), so row 7 of each locale
# status sheet (zh-cn / de-de) is populated with the handback results:
#   I = Latest Target File (new hyperlink to the .md on GitHub)
#   J = Latest Handback File (the generated .xlf)
#   K = Latest Handback DateTime
#   P = Error Detail (stale handback warning)

$rowData = @{
    "zh-cn" = @{
        XlfName = "f1f0eb99-5294-4e84-87c1-4da7274f6c19.751183d91aa34ec17002a77c08aadfc9c9d39cfc.zh-cn.xlf"
        HandbackDateTime = "2016-08-20 19:03:56"
    }
    "de-de" = @{
        XlfName = "f1f0eb99-5294-4e84-87c1-4da7274f6c19.751183d91aa34ec17002a77c08aadfc9c9d39cfc.de-de.xlf"
        HandbackDateTime = "2016-08-20 19:04:05"
    }
}

$mdName = "f1f0eb99-5294-4e84-87c1-4da7274f6c19.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b4d2500ea33572574f205b79182705b836ae5f7e/e2e/f1f0eb99-5294-4e84-87c1-4da7274f6c19.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8eb72e63cdcf1513722b5164085688dfb4d858cc/e2e/f1f0eb99-5294-4e84-87c1-4da7274f6c19.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b4d2500ea33572574f205b79182705b836ae5f7e/e2e/f1f0eb99-5294-4e84-87c1-4da7274f6c19.md."

foreach ($name in $rowData.Keys) {
    $ws = $wb.Worksheets.Item($name)
    $info = $rowData[$name]

    $ws.Range("I7").Value = $mdName
    $ws.Range("J7").Value = $info.XlfName
    $ws.Range("K7").Value = $info.HandbackDateTime
    $ws.Range("P7").Value = $errorDetail

    $ws.Hyperlinks.Add($ws.Range("I7"), $mdUrl, "", "", $mdName)
}
